# Listas sem duplicação de professores
# Replace list-like cell values (which represented duplicated/placeholder
# teacher-class lists) with a simple "-" placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @("B3", "F3", "B4", "D4", "F4", "B6", "D6", "F6", "B7", "D7", "F11", "F12", "F14", "B18", "C18", "D18", "B19", "D19", "B20", "B21", "D21")

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value = "-"
}
